$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B10: drop the 4th item from the list
$ws.Range("B10").Value = "ТС- ОБ-ПР1, ТС- ОБ-ПР2, ТС- ОБ-ПР3 "

# Update cells B18:B20 with new scenario labels
$ws.Range("B18").Value = "ТС-ДАТ-В1"
$ws.Range("B19").Value = "ТС-НАСТР1"
$ws.Range("B20").Value = "ТС-НАСТР2"

# Update the visible selection/scroll position of the sheet view
$ws.Activate()
$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
